# Updates cryptocurrency price/volume data in the active worksheet,
# reflecting refreshed rankings (including two pairs of rows whose
# coin identity swapped position: PEPE/Binance-PegBSC-USD and
# Aptos/Bittensor).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.659.79'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.64%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.425.90'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -6.46%  '
$ws.Range("E4").Value = '  -0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.86'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.75%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -6.38%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.422.81'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.45%  '
$ws.Range("E10").Value = '  -9.26%  '
$ws.Range("E11").Value = '  -1.85%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.34'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.73%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.345'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -8.39%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.55'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.94%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.860.97'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.86%  '
$ws.Range("E16").Value = '  -9.78%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.638.10'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.62%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.424.01'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -7.14%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -9.05%  '
$ws.Range("E20").Value = '  -8.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.10'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -7.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '314.35'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.89%  '
$ws.Range("E23").Value = '  -0.02%  '
$ws.Range("E24").Value = '  -1.98%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '62.94'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.48%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.556.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.93%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.22%  '
$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0941'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -13.53%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -10.73%  '
$ws.Range("E30").Value = '  -11.34%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '513.58'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -8.02%  '
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.54'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.19%  '
$ws.Range("E33").Value = '  -7.66%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.49%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.54'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.45%  '
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.56'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -13.30%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.69'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -11.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.372'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.13'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.46%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '140.59'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.69%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.39%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.01'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.20'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -10.79%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '138.52'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -12.78%  '
$ws.Range("E47").Value = '  -7.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '20.53'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -12.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0524'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.577'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -6.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0922'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.38%  '
